$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.829.21'
$ws.Range("E2").Value = '  -1.26%  '
$ws.Range("D3").Value = '1.634.05'
$ws.Range("E4").Value = '  -0.42%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '214.96'
$ws.Range("E5").Value = '  -0.50%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5016'
$ws.Range("E6").Value = '  -1.72%  '
$ws.Range("E7").Value = '  -0.42%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2566'
$ws.Range("E8").Value = '  -0.94%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06403'
$ws.Range("E9").Value = '  -0.11%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.60'
$ws.Range("E10").Value = '  -1.85%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07687'
$ws.Range("E11").Value = '  -1.64%  '
$ws.Range("B12").Value = 'Polkadot'
$ws.Range("C12").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.239'
$ws.Range("E12").Value = '  -1.23%  '
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = '1.632.59'
$ws.Range("E13").Value = '  -1.45%  '
$ws.Range("D14").Value = '1.859.66'
$ws.Range("E14").Value = '  -1.26%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.5424'
$ws.Range("E15").Value = '  -2.05%  '
$ws.Range("D16").Value = '0.0₅7929'
$ws.Range("E16").Value = '  -1.16%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '63.39'
$ws.Range("E17").Value = '  -1.21%  '
$ws.Range("D18").Value = '25.852.89'
$ws.Range("E18").Value = '  -1.29%  '
$ws.Range("E19").Value = '  -0.29%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '202.17'
$ws.Range("E20").Value = '  -3.51%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.324'
$ws.Range("E21").Value = '  -2.17%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.923'
$ws.Range("E22").Value = '  -1.41%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.969'
$ws.Range("E23").Value = '  -0.93%  '
$ws.Range("E24").Value = '  -0.35%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.909'
$ws.Range("E25").Value = '  +10.58%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '140.80'
$ws.Range("E26").Value = '  -1.97%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.1137'
$ws.Range("E27").Value = '  -2.74%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.67'
$ws.Range("E28").Value = '  -0.79%  '
$ws.Range("E29").Value = '  -4.10%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.239'
$ws.Range("E30").Value = '  -0.57%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.04978'
$ws.Range("E31").Value = '  -2.49%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.253'
$ws.Range("E32").Value = '  -2.96%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.173'
$ws.Range("E33").Value = '  -1.75%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.537'
$ws.Range("E34").Value = '  -1.67%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.362'
$ws.Range("E35").Value = '  -0.52%  '
$ws.Range("D36").Value = '1.165.88'
$ws.Range("E36").Value = '  +1.04%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.8917'
$ws.Range("E37").Value = '  -3.99%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.615'
$ws.Range("E38").Value = '  -4.91%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.5595'
$ws.Range("E39").Value = '  -1.99%  '
$ws.Range("E40").Value = '  -2.27%  '
$ws.Range("B41").Value = 'PaxDollar'
$ws.Range("C41").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.002'
$ws.Range("E41").Value = '  -0.37%  '
$ws.Range("B42").Value = 'mCoin'
$ws.Range("C42").Value = 'https://coinranking.com/coin/fzVgyjBcRc9+mcoin-mcoin'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.543'
$ws.Range("E42").Value = '  -0.72%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.674'
$ws.Range("E43").Value = '  +0.43%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.8066'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '99.20'
$ws.Range("E45").Value = '  -1.33%  '
$ws.Range("D46").Value = '1.771.69'
$ws.Range("E46").Value = '  -1.26%  '
$ws.Range("E47").Value = '  -0.70%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.4515'
$ws.Range("E48").Value = '  -0.73%  '
$ws.Range("E49").Value = '  +0.07%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '54.67'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05076'
$ws.Range("E51").Value = '  +0.55%  '
